# Generate Report for Handoff
# Moves the localization status from "In Translation" to "Ready for
# handoff" and refreshes the handoff timestamps that the report stamps
# when a handoff package is generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Status column everywhere it appears -------------------------------
$wsOverview.Range("E2").Value = $newStatus   # zh-cn status (Overview sheet)
$wsOverview.Range("F2").Value = $newStatus   # de-de status (Overview sheet)
$wsZhCn.Range("C2").Value     = $newStatus   # Status column, zh-cn sheet
$wsDeDe.Range("C2").Value     = $newStatus   # Status column, de-de sheet

# --- Handoff timestamps, bumped to the moment of this handoff ----------
$wsZhCn.Range("H2").Value     = "2016-08-16 20:54:40"   # Latest Handoff Datetime (zh-cn)
$wsDeDe.Range("H2").Value     = "2016-08-16 20:54:45"   # Latest Handoff Datetime (de-de)
$wsOverview.Range("G2").Value = "2016-08-16 20:54:45"   # Latest HO Xliff Generate Date

# --- Re-fit the Status / language columns for the new, longer text -----
# ColumnWidth is expressed in characters; Excel stores the column run on
# a pixel grid, so feed it the width (in characters) that lands on the
# pixel the refreshed report was rendered at.
$statusColumnWidth = 16.333333333333332

$wsOverview.Columns.Item(5).ColumnWidth = $statusColumnWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $statusColumnWidth   # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth     = $statusColumnWidth   # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = $statusColumnWidth   # column C (Status)
